$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header block: "Heavy" / "Hex Nut" / "Width" / "Across Flat"
# (was: (blank) / "Bolt" / "Head" / "Dia (in.)")
$ws.Range("K2").Value = "Heavy"
$ws.Range("K2").Font.Bold = $true

$ws.Range("K3").Value = "Hex Nut"
$ws.Range("K4").Value = "Width"
$ws.Range("K5").Value = "Across Flat"

# Page orientation changed to portrait
$ws.PageSetup.Orientation = 1

# Update selection to reflect where the user left off editing
$ws.Range("K6").Select() | Out-Null
